# Calificaciones Const y Var grupo C
# Fill in the "DE 0 A 5 CUMPLE?" grade column (D) for every student row,
# and flag the two students who did not present with a 0 and a
# "No presenta" note in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$noPresentaRows = @(10, 21)

for ($r = 3; $r -le 31; $r++) {
    if ($noPresentaRows -contains $r) {
        $ws.Cells.Item($r, 4).Value = 0
        $ws.Cells.Item($r, 5).Value = "No presenta"
    } else {
        $ws.Cells.Item($r, 4).Value = 5
    }
}

# Reflect the reviewer's final scroll position / zoom / selection on the sheet.
$excel.ActiveWindow.Zoom = 140
$ws.Range("F22").Select() | Out-Null
